$wb = $excel.ActiveWorkbook

# ==== 1. Create the new "2022-Q3" worksheet ====
# Duplicate "2022-Q2" (same column layout/header/styles) and drop it in right
# before "2022-Q2" so sheet order becomes 总计, 2022-Q3, 2022-Q2, 2022-Q1, ...
$sheetQ2Template = $wb.Worksheets.Item("2022-Q2")
$sheetQ2Template.Copy($sheetQ2Template, $null)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The template has 7 data+header rows; 2022-Q3 needs 12 (header + 11 funds) -
# stamp the last template row down to manufacture the extra styled rows.
for ($r = 8; $r -le 12; $r++) {
    $newSheet.Range("A7:H7").Copy($newSheet.Range("A" + $r + ":H" + $r))
}

# Columns B:G hold text (fund code/name/size/position/ratio/value) in this workbook,
# not numbers - force text so the "8.92"-style strings are not coerced to numeric.
$newSheet.Range("B2:G12").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "217002"
$newSheet.Range("C2").Value = "招商安泰平衡混合"
$newSheet.Range("D2").Value = "8.92"
$newSheet.Range("E2").Value = "49.42"
$newSheet.Range("F2").Value = "1.83"
$newSheet.Range("G2").Value = "0.1632"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "011160"
$newSheet.Range("C3").Value = "富国质量成长6个月持有期混合A"
$newSheet.Range("D3").Value = "3.70"
$newSheet.Range("E3").Value = "85.89"
$newSheet.Range("F3").Value = "3.19"
$newSheet.Range("G3").Value = "0.1180"
$newSheet.Range("H3").Value = 7

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "010122"
$newSheet.Range("C4").Value = "华泰柏瑞优势领航混合A"
$newSheet.Range("D4").Value = "3.95"
$newSheet.Range("E4").Value = "79.94"
$newSheet.Range("F4").Value = "2.41"
$newSheet.Range("G4").Value = "0.0952"
$newSheet.Range("H4").Value = 8

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "002317"
$newSheet.Range("C5").Value = "招商睿逸稳健配置混合"
$newSheet.Range("D5").Value = "6.22"
$newSheet.Range("E5").Value = "46.67"
$newSheet.Range("F5").Value = "1.42"
$newSheet.Range("G5").Value = "0.0883"
$newSheet.Range("H5").Value = 10

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "004671"
$newSheet.Range("C6").Value = "中融核心成长灵活配置混合"
$newSheet.Range("D6").Value = "1.10"
$newSheet.Range("E6").Value = "62.39"
$newSheet.Range("F6").Value = "2.90"
$newSheet.Range("G6").Value = "0.0319"
$newSheet.Range("H6").Value = 7

$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "008422"
$newSheet.Range("C7").Value = "中融研发创新混合A"
$newSheet.Range("D7").Value = "0.69"
$newSheet.Range("E7").Value = "61.89"
$newSheet.Range("F7").Value = "3.06"
$newSheet.Range("G7").Value = "0.0211"
$newSheet.Range("H7").Value = 8

$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "010008"
$newSheet.Range("C8").Value = "中融成长优选混合A"
$newSheet.Range("D8").Value = "0.58"
$newSheet.Range("E8").Value = "62.36"
$newSheet.Range("F8").Value = "2.95"
$newSheet.Range("G8").Value = "0.0171"
$newSheet.Range("H8").Value = 7

$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "010009"
$newSheet.Range("C9").Value = "中融成长优选混合C"
$newSheet.Range("D9").Value = "0.51"
$newSheet.Range("E9").Value = "62.36"
$newSheet.Range("F9").Value = "2.95"
$newSheet.Range("G9").Value = "0.0150"
$newSheet.Range("H9").Value = 7

$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "008423"
$newSheet.Range("C10").Value = "中融研发创新混合C"
$newSheet.Range("D10").Value = "0.40"
$newSheet.Range("E10").Value = "61.89"
$newSheet.Range("F10").Value = "3.06"
$newSheet.Range("G10").Value = "0.0122"
$newSheet.Range("H10").Value = 8

$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "010123"
$newSheet.Range("C11").Value = "华泰柏瑞优势领航混合C"
$newSheet.Range("D11").Value = "0.39"
$newSheet.Range("E11").Value = "79.94"
$newSheet.Range("F11").Value = "2.41"
$newSheet.Range("G11").Value = "0.0094"
$newSheet.Range("H11").Value = 8

$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "011161"
$newSheet.Range("C12").Value = "富国质量成长6个月持有期混合C"
$newSheet.Range("D12").Value = "0.14"
$newSheet.Range("E12").Value = "85.89"
$newSheet.Range("F12").Value = "3.19"
$newSheet.Range("G12").Value = "0.0045"
$newSheet.Range("H12").Value = 7

# Drop the temporary text format now that the values are committed as text, so no
# extra "@"-format style lingers on the range (matches a plain/General-formatted cell).
$newSheet.Range("B2:G12").Style = "Normal"

# ==== 2. Update the "总计" summary sheet: new 2022-Q3 row, older rows shift down ====
$totalSheet = $wb.Worksheets.Item("总计")

# Extend the styled index-column / row pattern down to the new row 6 (was A1:D5)
$totalSheet.Range("A5:D5").Copy($totalSheet.Range("A6:D6"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 11
$totalSheet.Range("D2").Value = 0.58

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 6
$totalSheet.Range("D3").Value = 2.86

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 7
$totalSheet.Range("D4").Value = 4.11

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 4
$totalSheet.Range("D5").Value = 2.34

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q3"
$totalSheet.Range("C6").Value = 2
$totalSheet.Range("D6").Value = 0.99

# ==== 3. Restore the original tab selection ====
# Copy()/Add() activates the new sheet; the source workbook had "2021-Q3" selected,
# so put the focus back there to leave bookViews/sheetViews untouched.
$wb.Worksheets.Item("2021-Q3").Activate()

